# Insert a new weekly price record as row 51, pushing all existing rows
# (51..122) down by one (52..123), per the "Fruta / hortaliza, semanal" update.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(51).Insert()

$ws.Range("A51").Value = 10
$ws.Range("B51").Value = "Vega Modelo de Temuco"
$ws.Range("C51").Value = "La Araucanía"
$ws.Range("D51").Value = 44803
$ws.Range("E51").Value = 9
$ws.Range("F51").Value = 100112035
$ws.Range("G51").Value = "Bruselas (repollito)"
$ws.Range("H51").Value = "Sin especificar"
$ws.Range("I51").Value = "Primera"
$ws.Range("J51").Value = 60
$ws.Range("K51").Value = 24000
$ws.Range("L51").Value = 24000
$ws.Range("M51").Value = 24000
$ws.Range("N51").Value = '$/malla 10 kilos'
$ws.Range("O51").Value = "Provincia de Quillota"
$ws.Range("P51").Value = 2400
$ws.Range("Q51").Value = 10
$ws.Range("R51").Value = "Hortaliza"
